$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row for the evaluation columns
$ws.Range("F1").Value = "Avaliação Matematica"
$ws.Range("G1").Value = "Avaliação Português"
$ws.Range("H1").Value = "Avaliação Ciências"

# Row 2 - João
$ws.Range("F2").Value = "Insuficiente"
$ws.Range("G2").Value = "Insuficiente"
$ws.Range("H2").Value = "Insuficiente"

# Row 3 - Maria
$ws.Range("F3").Value = "Suficiente"
$ws.Range("G3").Value = "Suficiente"
$ws.Range("H3").Value = "Suficiente"

# Row 4 - Pedro
$ws.Range("F4").Value = "Insuficiente"
$ws.Range("G4").Value = "Insuficiente"
$ws.Range("H4").Value = "Suficiente"

# Row 5 - Ana
$ws.Range("F5").Value = "Suficiente"
$ws.Range("G5").Value = "Suficiente"
$ws.Range("H5").Value = "Suficiente"
